$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.3300383186754654, 0.31678006652821755, -0.0046156228275238288),
    @(0.35249483409338667, 0.31824434937966328, -0.0046170387738013255),
    @(0.32627437462854975, 0.34063748112140146, -0.0046274803850593852),
    @(0.35249483409338667, 0.31824434937966328, -0.0046170387738013255),
    @(0.33471146151821679, 0.32311049088937494, -0.0046257674970918584),
    @(0.32627437462854975, 0.34063748112140146, -0.0046274803850593852),
    @(0.33471146151821679, 0.32311049088937494, -0.0046257674970918584),
    @(0.341493876083385, 0.32505916769252574, -0.0046268147292704478),
    @(0.32627437462854975, 0.34063748112140146, -0.0046274803850593852)
)

$startRow = 28
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
